$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
$ws.Range("D2").Value = "41.229.13"
$ws.Range("E2").Value = "  -0.18%  "
# Row 3
$ws.Range("D3").Value = "2.179.74"
$ws.Range("E3").Value = "  -1.70%  "
# Row 4
$ws.Range("E4").Value = "  -0.08%  "
# Row 5
Set-TextCell "D5" "255.98"
$ws.Range("E5").Value = "  +5.48%  "
# Row 6
Set-TextCell "D6" "0.625"
$ws.Range("E6").Value = "  -0.43%  "
# Row 7
Set-TextCell "D7" "67.74"
$ws.Range("E7").Value = "  -2.74%  "
# Row 8
$ws.Range("E8").Value = "  -0.06%  "
# Row 9
$ws.Range("E9").Value = "  +1.04%  "
# Row 10
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell "D10" "58.94"
$ws.Range("E10").Value = "  +1.66%  "
# Row 11
$ws.Range("B11").Value = "Avalanche"
$ws.Range("C11").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell "D11" "36.92"
$ws.Range("E11").Value = "  -6.12%  "
# Row 12
Set-TextCell "D12" "0.0930"
$ws.Range("E12").Value = "  -2.32%  "
# Row 13
Set-TextCell "D13" "6.99"
$ws.Range("E13").Value = "  +3.39%  "
# Row 14
$ws.Range("E14").Value = "  -0.70%  "
# Row 15
$ws.Range("D15").Value = "2.507.26"
$ws.Range("E15").Value = "  -1.53%  "
# Row 16
Set-TextCell "D16" "0.864"
$ws.Range("E16").Value = "  +2.63%  "
# Row 17
Set-TextCell "D17" "14.34"
$ws.Range("E17").Value = "  -3.25%  "
# Row 18
$ws.Range("D18").Value = "2.170.16"
$ws.Range("E18").Value = "  -2.08%  "
# Row 19
$ws.Range("D19").Value = "41.085.00"
$ws.Range("E19").Value = "  -0.39%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0952"
$ws.Range("E20").Value = "  +0.03%  "
# Row 21
$ws.Range("E21").Value = "  +0.59%  "
# Row 22
Set-TextCell "D22" "71.65"
$ws.Range("E22").Value = "  -0.70%  "
# Row 23
Set-TextCell "D23" "231.59"
$ws.Range("E23").Value = "  -0.15%  "
# Row 24
$ws.Range("E24").Value = "  -1.62%  "
# Row 25
Set-TextCell "D25" "3.91"
$ws.Range("E25").Value = "  +8.83%  "
# Row 26
$ws.Range("E26").Value = "  +19.15%  "
# Row 27
Set-TextCell "D27" "1.00"
$ws.Range("E27").Value = "  -0.01%  "
# Row 28
Set-TextCell "D28" "2.51"
$ws.Range("E28").Value = "  +4.07%  "
# Row 29
$ws.Range("E29").Value = "  -0.60%  "
# Row 30
Set-TextCell "D30" "168.73"
$ws.Range("E30").Value = "  -2.02%  "
# Row 31
Set-TextCell "D31" "20.56"
$ws.Range("E31").Value = "  +0.56%  "
# Row 32
Set-TextCell "D32" "0.116"
$ws.Range("E32").Value = "  -2.75%  "
# Row 33
Set-TextCell "D33" "0.0749"
$ws.Range("E33").Value = "  +4.36%  "
# Row 34
$ws.Range("E34").Value = "  -0.99%  "
# Row 35
Set-TextCell "D35" "5.44"
$ws.Range("E35").Value = "  +3.99%  "
# Row 36
Set-TextCell "D36" "26.29"
$ws.Range("E36").Value = "  +9.98%  "
# Row 37
Set-TextCell "D37" "4.12"
$ws.Range("E37").Value = "  +4.79%  "
# Row 38
Set-TextCell "D38" "4.59"
$ws.Range("E38").Value = "  -0.50%  "
# Row 39
Set-TextCell "D39" "0.0295"
$ws.Range("E39").Value = "  +5.87%  "
# Row 40
$ws.Range("E40").Value = "  -3.97%  "
# Row 41
Set-TextCell "D41" "12.26"
$ws.Range("E41").Value = "  +13.60%  "
# Row 42
Set-TextCell "D42" "5.63"
$ws.Range("E42").Value = "  -3.79%  "
# Row 43
Set-TextCell "D43" "63.32"
$ws.Range("E43").Value = "  -2.91%  "
# Row 44
Set-TextCell "D44" "5.02"
$ws.Range("E44").Value = "  -0.85%  "
# Row 45
$ws.Range("E45").Value = "  -1.42%  "
# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell "D46" "8.56"
$ws.Range("E46").Value = "  -1.70%  "
# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell "D47" "0.101"
$ws.Range("E47").Value = "  +0.50%  "
# Row 48
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextCell "D48" "1.19"
$ws.Range("E48").Value = "  +8.04%  "
# Row 49
$ws.Range("B49").Value = "BinanceUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D49" "1.00"
$ws.Range("E49").Value = "  +0.05%  "
# Row 50
$ws.Range("E50").Value = "  -0.32%  "
# Row 51
$ws.Range("E51").Value = "  -6.83%  "
